$d = $word.ActiveDocument
$d.Content.Find.Execute("In due settimane", $true, $false, $false, $false, $false, $true, 1, $false, "In tre settimane", 2)
